# Auto-generated Excel COM-interop script
# Applies updated market-price / profit figures to several Leve tables
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR worksheets.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 2307.2942
$ws.Range("I80").Value = 1420.8422
$ws.Range("K80").Value = 4262.5266
$ws.Range("M80").Value = -3264.5266
$ws.Range("H83").Value = 2307.2942
$ws.Range("I83").Value = 1420.8422
$ws.Range("K83").Value = 12787.5798
$ws.Range("M83").Value = -7795.5798
$ws.Range("H108").Value = 46500
$ws.Range("J108").Value = 46500
$ws.Range("L108").Value = 46500
$ws.Range("N108").Value = -54180
$ws.Range("H110").Value = 60000
$ws.Range("J110").Value = 60000
$ws.Range("L110").Value = 60000
$ws.Range("N110").Value = -68180
$ws.Range("H116").Value = 2873.375
$ws.Range("I116").Value = 2331.3333
$ws.Range("K116").Value = 2331.3333
$ws.Range("M116").Value = 1110.6667
$ws.Range("H137").Value = 5811.8423
$ws.Range("I137").Value = 7048.8335
$ws.Range("K137").Value = 21146.5005
$ws.Range("M137").Value = -18596.5005
$ws.Range("H140").Value = 37665.668
$ws.Range("I140").Value = 32999
$ws.Range("J140").Value = 39999
$ws.Range("K140").Value = 32999
$ws.Range("L140").Value = 39999
$ws.Range("M140").Value = -27819
$ws.Range("N140").Value = -50359
$ws.Range("H141").Value = 3190.6667
$ws.Range("I141").Value = 3190.6667
$ws.Range("K141").Value = 9572.000100000001
$ws.Range("M141").Value = -4392.000100000001

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 877
$ws.Range("I19").Value = 877
$ws.Range("K19").Value = 877
$ws.Range("M19").Value = -648
$ws.Range("H122").Value = 4852.933
$ws.Range("I122").Value = 3256.4285
$ws.Range("K122").Value = 9769.2855
$ws.Range("M122").Value = -7319.2855
$ws.Range("H130").Value = 154832.16
$ws.Range("J130").Value = 154832.16
$ws.Range("L130").Value = 154832.16
$ws.Range("N130").Value = -164872.16
$ws.Range("H132").Value = 2311.561
$ws.Range("I132").Value = 1494.5938
$ws.Range("K132").Value = 4483.7814
$ws.Range("M132").Value = -1953.7814
$ws.Range("H138").Value = 66618.125
$ws.Range("J138").Value = 66618.125
$ws.Range("L138").Value = 66618.125
$ws.Range("N138").Value = -76898.125

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3216.8333
$ws.Range("I99").Value = 3247.75
$ws.Range("J99").Value = 3155
$ws.Range("K99").Value = 3247.75
$ws.Range("L99").Value = 3155
$ws.Range("M99").Value = -1749.75
$ws.Range("N99").Value = -6151
$ws.Range("H128").Value = 4999
$ws.Range("I128").Value = 4999
$ws.Range("K128").Value = 14997
$ws.Range("M128").Value = -12507

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2496.3044
$ws.Range("I31").Value = 1599.2
$ws.Range("J31").Value = 5350.727
$ws.Range("K31").Value = 1599.2
$ws.Range("L31").Value = 5350.727
$ws.Range("M31").Value = -1304.2
$ws.Range("N31").Value = -5940.727
$ws.Range("H34").Value = 2496.3044
$ws.Range("I34").Value = 1599.2
$ws.Range("J34").Value = 5350.727
$ws.Range("K34").Value = 1599.2
$ws.Range("L34").Value = 5350.727
$ws.Range("M34").Value = -1397.2
$ws.Range("N34").Value = -5754.727
$ws.Range("H120").Value = 1712666.4
$ws.Range("J120").Value = 1712666.4
$ws.Range("L120").Value = 1712666.4
$ws.Range("N120").Value = -1719924.4
$ws.Range("H134").Value = 4233.685
$ws.Range("I134").Value = 2923.3225
$ws.Range("K134").Value = 8769.967500000001
$ws.Range("M134").Value = -6234.967500000001

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 947
$ws.Range("J122").Value = 1034.3
$ws.Range("L122").Value = 9308.699999999999
$ws.Range("N122").Value = -14208.7
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("H134").Value = 10671.667
$ws.Range("I134").Value = 1676.6666
$ws.Range("K134").Value = 5029.9998
$ws.Range("M134").Value = 40.0002000000004

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 6954667
$ws.Range("I11").Value = 10240500
$ws.Range("J11").Value = 383000
$ws.Range("K11").Value = 10240500
$ws.Range("L11").Value = 383000
$ws.Range("M11").Value = -10240361
$ws.Range("N11").Value = -383278
$ws.Range("H70").Value = 46752.242
$ws.Range("J70").Value = 10254
$ws.Range("L70").Value = 10254
$ws.Range("N70").Value = -10794
$ws.Range("H73").Value = 46752.242
$ws.Range("J73").Value = 10254
$ws.Range("L73").Value = 10254
$ws.Range("N73").Value = -12126
$ws.Range("H126").Value = 5116.222
$ws.Range("J126").Value = 5505.75
$ws.Range("L126").Value = 16517.25
$ws.Range("N126").Value = -21457.25
$ws.Range("H131").Value = 91353.39999999999
$ws.Range("J131").Value = 91353.39999999999
$ws.Range("L131").Value = 91353.39999999999
$ws.Range("N131").Value = -101433.4

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("H88").Value = 59995
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 59995
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 59995
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -60851
$ws.Range("H91").Value = 59995
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 59995
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 59995
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -62959
$ws.Range("H136").Value = 4812.2666
$ws.Range("I136").Value = 3016.7693
$ws.Range("J136").Value = 6185.294
$ws.Range("K136").Value = 9050.3079
$ws.Range("L136").Value = 18555.882
$ws.Range("M136").Value = -6500.3079
$ws.Range("N136").Value = -23655.882

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2741.8838
$ws.Range("I136").Value = 2298.9092
$ws.Range("J136").Value = 4203.7
$ws.Range("K136").Value = 6896.7276
$ws.Range("L136").Value = 12611.1
$ws.Range("M136").Value = -4346.7276
$ws.Range("N136").Value = -17711.1
$ws.Range("H138").Value = 75430.625
$ws.Range("J138").Value = 75430.625
$ws.Range("L138").Value = 75430.625
$ws.Range("N138").Value = -85710.625
